$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 799.31915
$ws.Range("J17").Value = 817.06665
$ws.Range("L17").Value = 2451.19995
$ws.Range("N17").Value = -2787.19995
$ws.Range("H38").Value = 2446.2222
$ws.Range("I38").Value = 335.66666
$ws.Range("J38").Value = 3501.5
$ws.Range("K38").Value = 1006.99998
$ws.Range("L38").Value = 10504.5
$ws.Range("M38").Value = -634.9999799999999
$ws.Range("N38").Value = -11248.5
$ws.Range("H137").Value = 142859400
$ws.Range("I137").Value = 200002350
$ws.Range("J137").Value = 2000
$ws.Range("K137").Value = 600007050
$ws.Range("L137").Value = 6000
$ws.Range("M137").Value = -600004500
$ws.Range("N137").Value = -11100
$ws.Range("H138").Value = 6739509
$ws.Range("I138").Value = 2151515.2
$ws.Range("J138").Value = 7814820.5
$ws.Range("K138").Value = 6454545.600000001
$ws.Range("L138").Value = 23444461.5
$ws.Range("M138").Value = -6449405.600000001
$ws.Range("N138").Value = -23454741.5

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 16839.463
$ws.Range("I32").Value = 1551.6964
$ws.Range("J32").Value = 82694.46000000001
$ws.Range("K32").Value = 1551.6964
$ws.Range("L32").Value = 82694.46000000001
$ws.Range("M32").Value = -1264.6964
$ws.Range("N32").Value = -83268.46000000001
$ws.Range("H61").Value = 2521.2
$ws.Range("I61").Value = 1901.5
$ws.Range("J61").Value = 5000
$ws.Range("K61").Value = 1901.5
$ws.Range("L61").Value = 5000
$ws.Range("M61").Value = -1689.5
$ws.Range("N61").Value = -5424
$ws.Range("H74").Value = 9007.177
$ws.Range("I74").Value = 1941.4667
$ws.Range("J74").Value = 62000
$ws.Range("K74").Value = 1941.4667
$ws.Range("L74").Value = 62000
$ws.Range("M74").Value = -1067.4667
$ws.Range("N74").Value = -63748
$ws.Range("H77").Value = 9007.177
$ws.Range("I77").Value = 1941.4667
$ws.Range("J77").Value = 62000
$ws.Range("K77").Value = 9707.333500000001
$ws.Range("L77").Value = 310000
$ws.Range("M77").Value = -5339.333500000001
$ws.Range("N77").Value = -318736
$ws.Range("H102").Value = 2980.2
$ws.Range("I102").Value = 3500.25
$ws.Range("K102").Value = 3500.25
$ws.Range("M102").Value = -1878.25
$ws.Range("H123").Value = 32942.8
$ws.Range("J123").Value = 32942.8
$ws.Range("L123").Value = 32942.8
$ws.Range("N123").Value = -42742.8
$ws.Range("H132").Value = 2967.2188
$ws.Range("I132").Value = 2475.9167
$ws.Range("J132").Value = 4441.125
$ws.Range("K132").Value = 7427.750100000001
$ws.Range("L132").Value = 13323.375
$ws.Range("M132").Value = -4897.750100000001
$ws.Range("N132").Value = -18383.375
$ws.Range("H136").Value = 2521.2
$ws.Range("I136").Value = 1901.5
$ws.Range("J136").Value = 5000
$ws.Range("K136").Value = 5704.5
$ws.Range("L136").Value = 15000
$ws.Range("M136").Value = -3154.5
$ws.Range("N136").Value = -20100
$ws.Range("H139").Value = 59000
$ws.Range("J139").Value = 59000
$ws.Range("L139").Value = 59000
$ws.Range("N139").Value = -69280

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H112").Value = 0
$ws.Range("J112").Value = 0
$ws.Range("L112").Value = 0
$ws.Range("N112").ClearContents()
$ws.Range("H133").Value = 53000
$ws.Range("J133").Value = 53000
$ws.Range("L133").Value = 53000
$ws.Range("N133").Value = -63120
$ws.Range("H134").Value = 3215.6
$ws.Range("I134").Value = 2048.2964
$ws.Range("K134").Value = 6144.889200000001
$ws.Range("M134").Value = -3609.889200000001

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 1338.9688
$ws.Range("I31").Value = 1062.3549
$ws.Range("J31").Value = 9914
$ws.Range("K31").Value = 1062.3549
$ws.Range("L31").Value = 9914
$ws.Range("M31").Value = -767.3549
$ws.Range("N31").Value = -10504
$ws.Range("H34").Value = 1338.9688
$ws.Range("I34").Value = 1062.3549
$ws.Range("J34").Value = 9914
$ws.Range("K34").Value = 1062.3549
$ws.Range("L34").Value = 9914
$ws.Range("M34").Value = -860.3549
$ws.Range("N34").Value = -10318
$ws.Range("H86").Value = 33335348
$ws.Range("I86").Value = 62501664
$ws.Range("J86").Value = 2415.1428
$ws.Range("K86").Value = 62501664
$ws.Range("L86").Value = 2415.1428
$ws.Range("M86").Value = -62500541
$ws.Range("N86").Value = -4661.1428
$ws.Range("H89").Value = 33335348
$ws.Range("I89").Value = 62501664
$ws.Range("J89").Value = 2415.1428
$ws.Range("K89").Value = 312508320
$ws.Range("L89").Value = 12075.714
$ws.Range("M89").Value = -312502704
$ws.Range("N89").Value = -23307.714
$ws.Range("H132").Value = 4360.923
$ws.Range("I132").Value = 3230
$ws.Range("J132").Value = 5067.75
$ws.Range("K132").Value = 9690
$ws.Range("L132").Value = 15203.25
$ws.Range("M132").Value = -7160
$ws.Range("N132").Value = -20263.25
$ws.Range("H135").Value = 39648.277
$ws.Range("J135").Value = 39648.277
$ws.Range("L135").Value = 39648.277
$ws.Range("N135").Value = -49788.277

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H122").Value = 847.44446
$ws.Range("I122").Value = 652
$ws.Range("J122").Value = 1003.8
$ws.Range("K122").Value = 5868
$ws.Range("L122").Value = 9034.199999999999
$ws.Range("M122").Value = -3418
$ws.Range("N122").Value = -13934.2
$ws.Range("H131").Value = 2676.065
$ws.Range("I131").Value = 308
$ws.Range("J131").Value = 2840.514
$ws.Range("K131").Value = 924
$ws.Range("L131").Value = 8521.542000000001
$ws.Range("M131").Value = 4116
$ws.Range("N131").Value = -18601.542
$ws.Range("H132").Value = 1141.7646
$ws.Range("J132").Value = 1136.1818
$ws.Range("L132").Value = 10225.6362
$ws.Range("N132").Value = -15285.6362

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H107").Value = 1281.3334
$ws.Range("H132").Value = 4865
$ws.Range("I132").Value = 4514.75
$ws.Range("J132").Value = 5332
$ws.Range("K132").Value = 13544.25
$ws.Range("L132").Value = 15996
$ws.Range("M132").Value = -11014.25
$ws.Range("N132").Value = -21056
$ws.Range("H138").Value = 54533.332
$ws.Range("J138").Value = 54533.332
$ws.Range("L138").Value = 54533.332
$ws.Range("N138").Value = -64813.332

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H17").Value = 2601.6
$ws.Range("I17").Value = 2004
$ws.Range("J17").Value = 3000
$ws.Range("K17").Value = 2004
$ws.Range("L17").Value = 3000
$ws.Range("M17").Value = -1834
$ws.Range("N17").Value = -3340
$ws.Range("H22").Value = 18138.834
$ws.Range("I22").Value = 0
$ws.Range("J22").Value = 18138.834
$ws.Range("K22").Value = 0
$ws.Range("L22").Value = 18138.834
$ws.Range("M22").ClearContents()
$ws.Range("N22").Value = -18728.834
$ws.Range("H27").Value = 18138.834
$ws.Range("I27").Value = 0
$ws.Range("J27").Value = 18138.834
$ws.Range("K27").Value = 0
$ws.Range("L27").Value = 18138.834
$ws.Range("M27").ClearContents()
$ws.Range("N27").Value = -18352.834
$ws.Range("H46").Value = 1695.25
$ws.Range("I46").Value = 450.5
$ws.Range("J46").Value = 2940
$ws.Range("K46").Value = 450.5
$ws.Range("L46").Value = 2940
$ws.Range("M46").Value = -262.5
$ws.Range("N46").Value = -3316
$ws.Range("H132").Value = 4610.722
$ws.Range("I132").Value = 3496
$ws.Range("J132").Value = 5725.4443
$ws.Range("K132").Value = 10488
$ws.Range("L132").Value = 17176.3329
$ws.Range("M132").Value = -7958
$ws.Range("N132").Value = -22236.3329

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 55560224
$ws.Range("I132").Value = 71433144
$ws.Range("J132").Value = 4999
$ws.Range("K132").Value = 214299432
$ws.Range("L132").Value = 14997
$ws.Range("M132").Value = -214296902
$ws.Range("N132").Value = -20057

Write-Output "Applied 202 cell updates across 8 sheets."
